$wb = $excel.ActiveWorkbook

# --- Scenarios sheet -------------------------------------------------
$scenarios = $wb.Worksheets.Item("Scenarios")

# Row 3: Execute flag flips from Y to N (test class / case stay the same)
$scenarios.Range("A3").Value = "N"

# Row 4: new scenario - VisitTests / bookVisit, not executed
$scenarios.Range("A4").Value = "N"
$scenarios.Range("B4").Value = "patient.tests.VisitTests"
$scenarios.Range("D4").Value = "bookVisit"

# Match the pasted-in formatting on the new Test Case cell (wrap text,
# small SimSun font, white fill) as seen in the authored workbook.
$scenarios.Range("D4").Font.Name = "SimSun"
$scenarios.Range("D4").Font.Size = 9
$scenarios.Range("D4").Font.Color = 0
$scenarios.Range("D4").Interior.Color = 16777215
$scenarios.Range("D4").WrapText = $true

$scenarios.Range("A4").Select() | Out-Null

# --- Parameters sheet -------------------------------------------------
$parameters = $wb.Worksheets.Item("Parameters")

# New Sauce Labs parameters (row 8 left blank, matching the authored sheet)
$parameters.Range("A9").Value = "USERNAME"
$parameters.Range("B9").Value = "qaheal"

$parameters.Range("A10").Value = "ACCESS_KEY"
$parameters.Range("B10").Value = "e14bb2d7-155b-4775-8978-9365c5b22012"

$parameters.Range("A11").Value = "RemoteURL"
$parameters.Range("B11").Value = "ondemand.saucelabs.com:443/wd/hub"

$parameters.Range("B27").Select() | Out-Null
$parameters.Activate() | Out-Null
